$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @{
    2 = @(1.02, 1.022019712320756, 1.026044126995927, 1.025640728492428, 1.020460602463755, 1.02829654538974, 1.027207903667633, 1.02886792516688, 1.028465706534977, 1.023300821707284, 1.012991519222763)
    3 = @(1.02, 1.023525159687274, 1.027135790720846, 1.027094346827127, 1.022613375645724, 1.028632487468846, 1.028348331590699, 1.029766249445237, 1.029724917874062, 1.025256149313745, 1.013380779767709)
    4 = @(1.02, 1.024496528411916, 1.027839584577629, 1.02803258601876, 1.024003013151366, 1.028846947251065, 1.029083191217741, 1.030344362821838, 1.030536867888398, 1.026517700457992, 1.013631308607835)
    5 = @(1.02, 1.024904243914455, 1.02813484805636, 1.028426470233213, 1.024586437448701, 1.028936410949276, 1.029391399812052, 1.030586651720396, 1.030877540812095, 1.027047197196588, 1.013736311182931)
    6 = @(1.02, 1.024972663493331, 1.028184388454548, 1.028492573095892, 1.024684351915998, 1.028951391623156, 1.029443107066613, 1.0306272893202, 1.030934702243235, 1.027136052236284, 1.013753922917862)
    7 = @(1.02, 1.024501978855803, 1.027843532294856, 1.028037851271183, 1.024010811926535, 1.028848145396564, 1.029087312357342, 1.030347603234929, 1.030541422598974, 1.026524778967771, 1.013632712909429)
    8 = @(1.02, 1.022529063278486, 1.026413599783025, 1.026132478557351, 1.021188849844538, 1.028410684549718, 1.027593958704788, 1.029172177204973, 1.028891857793171, 1.023962407673111, 1.013123352820663)
    9 = @(1.02, 1.019030882606891, 1.023873732546306, 1.022756488051373, 1.016189498674567, 1.027617335330784, 1.024938514177074, 1.027076381908699, 1.025962881862427, 1.019418084671903, 1.012215323088286)
    10 = @(1.02, 1.016683443684286, 1.022166462154113, 1.02049271317272, 1.012837202036167, 1.027073110687905, 1.023151505763754, 1.025662225530792, 1.02399464807015, 1.0163676525107, 1.011602728271609)
    11 = @(1.02, 1.015663172451184, 1.021423759675263, 1.019509208130236, 1.011380695831438, 1.026833773690606, 1.022373612639091, 1.02504575445808, 1.023138542060665, 1.015041532804862, 1.011335705115729)
    12 = @(1.02, 1.015283610317027, 1.021147360090361, 1.019143384968987, 1.010838914387732, 1.026744315636263, 1.022084039760167, 1.024816139834359, 1.022819955636848, 1.014548136679512, 1.011236251621996)
    13 = @(1.02, 1.015365054530683, 1.021206672733428, 1.019221878376515, 1.010955163653602, 1.026763529993957, 1.022146182767681, 1.024865421597525, 1.022888320430688, 1.014654009173491, 1.011257596980033)
    14 = @(1.02, 1.015631809807815, 1.021400923201811, 1.019478979466256, 1.011335927883577, 1.026826390459727, 1.022349689349181, 1.025026787342462, 1.023112219726795, 1.015000765327182, 1.011327489773637)
    15 = @(1.02, 1.015796088178578, 1.021520537231043, 1.019637320530172, 1.011570426389133, 1.026865046850115, 1.02247499280258, 1.025126126374231, 1.023250092846084, 1.015214304309073, 1.011370517239277)
    16 = @(1.02, 1.016751074090923, 1.022215679518955, 1.020557915004054, 1.01293375909658, 1.027088916724025, 1.023203044445659, 1.025703050860795, 1.02405138287931, 1.016455549657078, 1.011620412184869)
    17 = @(1.02, 1.017349079793266, 1.022650795405348, 1.021134493168594, 1.013787599507053, 1.027228355133294, 1.023658624280745, 1.026063827753919, 1.024552972049033, 1.017232723197369, 1.011776689220401)
    18 = @(1.02, 1.01769751935318, 1.022904259578748, 1.021470486330652, 1.014285155482654, 1.027309332015868, 1.023923960856973, 1.026273864747786, 1.024845170119115, 1.017685529369995, 1.011867672890107)
    19 = @(1.02, 1.017816266427328, 1.022990628359056, 1.021584998270194, 1.014454729540458, 1.027336882919136, 1.02401436705323, 1.026345414687328, 1.024944739656563, 1.017839839515723, 1.011898667240174)
    20 = @(1.02, 1.017284957513651, 1.022604145952323, 1.021072664480481, 1.013696039813315, 1.027213431479641, 1.023609785862365, 1.026025161046408, 1.024499194680789, 1.017149392354503, 1.011759939794268)
    21 = @(1.02, 1.015553273378513, 1.021343735900279, 1.019403283669496, 1.011223823796312, 1.026807895046377, 1.022289779184572, 1.024979286611739, 1.023046303367031, 1.014898677021519, 1.011306915551289)
    22 = @(1.02, 1.014461085519589, 1.020548213266932, 1.018350745658479, 1.009664973223068, 1.026549689781825, 1.021456194982751, 1.024318055340386, 1.022129390001735, 1.01347882763332, 1.011020521700969)
    23 = @(1.02, 1.01504040314074, 1.020970227452236, 1.018908998441379, 1.010491782598886, 1.026686876692442, 1.021898443076818, 1.024668935431535, 1.022615791777561, 1.01423197414625, 1.011172493657097)
    24 = @(1.02, 1.017313932745712, 1.022625225856664, 1.021100603176889, 1.013737413163785, 1.027220175933844, 1.023631855063169, 1.026042634097105, 1.024523495503113, 1.01718704753832, 1.01176750866449)
    25 = @(1.02, 1.019937889769201, 1.024532785436878, 1.023631522979788, 1.017485262155884, 1.027825120977646, 1.025627911419862, 1.027621150544848, 1.026722791892952, 1.020596484193257, 1.012451332055029)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $ws.Cells.Item([int]$row, $colIndexes[$i]).Value = $vals[$i]
    }
}
